$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns A and B from 12 to 20
# (Excel's ColumnWidth property adds ~0.8333 padding when stored in the
# sheet XML "width" attribute, so we back-solve for the ColumnWidth value
# that serializes to exactly width="20".)
$ws.Columns.Item(1).ColumnWidth = 19.16666666666667
$ws.Columns.Item(2).ColumnWidth = 19.16666666666667

# Update row 2 values
$ws.Range("A2").Value = "mongodb_qa_micro01"
$ws.Range("B2").Value = "mongodb_qa_micro01"
$ws.Range("C2").Value = "10.181.34.51"
$ws.Range("D2").Value = "20.0.3"
$ws.Range("E2").Value = "860.ol8"
$ws.Range("F2").Value = "20.0.3"
$ws.Range("G2").Value = "860.ol8"
$ws.Range("H2").Value = "1.2.0.1253"
$ws.Range("I2").Value = "1.2.0.1253"
$ws.Range("K2").Value = "active"
$ws.Range("L2").Value = "active"
$ws.Range("M2").Value = "active"
